# "second version of erd"
# Reworks the Transaction / Loan-entity / Article area of the ERD on Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Drop the old "Loan entity" label that used to sit in B17 ---
$ws.Cells.Item(17,2).ClearContents()

# --- Introduce the new label text (order chosen to match the authoring
#     session's shared-string table growth) ---
$ws.Cells.Item(17,7).Value = "(deposit/ food/ items/ loan/ registration/ camp)"
$ws.Cells.Item(26,6).Value = "articleNr (pk)"
$ws.Cells.Item(25,6).Value = "Article"
$ws.Cells.Item(28,6).Value = "Condition (working / broken)"
$ws.Cells.Item(11,6).Value = "TransactionID (pk)"
$ws.Cells.Item(21,6).Value = "articleNr (pk) (fk)"
$ws.Cells.Item(27,6).Value = "articleType"
$ws.Cells.Item(22,6).Value = "TransactionID (pk) (fk)"
$ws.Cells.Item(20,6).Value = "Loan entity (intersection)"

# --- Shift the rest of the "Transaction" attribute list down by one row
#     underneath the newly inserted TransactionID (pk) row ---
$ws.Cells.Item(16,6).Value = "CurrentBalance"
$ws.Cells.Item(15,6).Value = "Amount"
$ws.Cells.Item(14,6).Value = "AccountID"
$ws.Cells.Item(13,6).Value = "Time"
$ws.Cells.Item(12,6).Value = "Date"

# F16 no longer carries the "Neutral" highlight (that moves to the new F17 row)
$ws.Cells.Item(16,6).Style = "Normal"

# --- Row 17: the "TransactionType" header moves into F17/G17 ---
$ws.Cells.Item(17,6).Value = "TransactionType"
$ws.Cells.Item(17,6).Style = "Neutral"

# --- Style the new section headers ---
$ws.Cells.Item(20,6).Style = "Accent1"
$ws.Cells.Item(25,6).Style = "Good"

# --- Widen column G to fit the new description text ---
$ws.Columns.Item(7).ColumnWidth = 32.6

# --- Move the selection cursor to match the saved cursor position ---
$ws.Activate()
$ws.Range("H27").Select()
